$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.195161700248718
$ws.Range("B1").Value = 2.176342248916626
$ws.Range("C1").Value = -1
$ws.Range("D1").Value = 2.215459108352661
$ws.Range("E1").Value = 1.213807821273804
